$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells that would otherwise be auto-converted to numbers stay as text,
# matching the original inline-string storage of these price cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the refreshed crypto feed
$ws.Range("D2").Value = "42.019.29"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "2.257.84"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "298.52"
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("D6").Value = "93.82"
$ws.Range("E6").Value = "  -6.97%  "
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").Value = "32.92"
$ws.Range("E10").Value = "  -5.60%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "47.44"
$ws.Range("E12").Value = "  -8.95%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "2.608.75"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").Value = "15.19"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").Value = "2.260.11"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "0.773"
$ws.Range("E18").Value = "  -4.27%  "
$ws.Range("D19").Value = "42.014.47"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("D22").Value = "11.36"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "66.52"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "233.20"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  -4.99%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").Value = "23.67"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "167.63"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("D31").Value = "33.56"
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "9.03"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "4.92"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -5.89%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("D37").Value = "0.0692"
$ws.Range("E37").Value = "  -4.90%  "
$ws.Range("D38").Value = "2.78"
$ws.Range("E38").Value = "  -6.23%  "
$ws.Range("E39").Value = "  -8.45%  "
$ws.Range("D40").Value = "0.0988"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("E42").Value = "  -8.58%  "
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").Value = "1.947.26"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "17.34"
$ws.Range("E46").Value = "  -7.76%  "
$ws.Range("E47").Value = "  -7.52%  "
$ws.Range("E48").Value = "  -4.85%  "
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("D50").Value = "2.482.63"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").Value = "51.94"
$ws.Range("E51").Value = "  -7.78%  "
